# Auto-generated Excel COM-interop script
# Applies numeric cell updates across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the authoritative diff of Phantom_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 656.125
$ws.Range("I9").Value = 818.1667
$ws.Range("J9").Value = 170
$ws.Range("K9").Value = 818.1667
$ws.Range("L9").Value = 170
$ws.Range("M9").Value = -649.1667
$ws.Range("N9").Value = -508
$ws.Range("H28").Value = 2422.5908
$ws.Range("I28").Value = 2301.4
$ws.Range("K28").Value = 2301.4
$ws.Range("M28").Value = -1816.4
$ws.Range("H40").Value = 1450.1666
$ws.Range("I40").Value = 1450.1666
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1450.1666
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1275.1666
$ws.Range("H51").Value = 16299.8
$ws.Range("I51").Value = 18250
$ws.Range("J51").Value = 14999.667
$ws.Range("K51").Value = 18250
$ws.Range("L51").Value = 14999.667
$ws.Range("M51").Value = -17766
$ws.Range("N51").Value = -15967.667
$ws.Range("H61").Value = 799.6667
$ws.Range("I61").Value = 799.6667
$ws.Range("K61").Value = 2399.0001
$ws.Range("M61").Value = -2227.0001
$ws.Range("H69").Value = 17628.455
$ws.Range("I69").Value = 14637.667
$ws.Range("K69").Value = 43913.001
$ws.Range("M69").Value = -43039.001
$ws.Range("H72").Value = 17628.455
$ws.Range("I72").Value = 14637.667
$ws.Range("K72").Value = 131739.003
$ws.Range("M72").Value = -127371.003
$ws.Range("H98").Value = 1084.25
$ws.Range("I98").Value = 1167.5555
$ws.Range("J98").Value = 334.5
$ws.Range("K98").Value = 1167.5555
$ws.Range("L98").Value = 334.5
$ws.Range("M98").Value = 330.4445000000001
$ws.Range("N98").Value = -3330.5
$ws.Range("H100").Value = 1568.4546
$ws.Range("I100").Value = 925.3
$ws.Range("J100").Value = 8000
$ws.Range("K100").Value = 925.3
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = -384.3
$ws.Range("N100").Value = -9082
$ws.Range("H122").Value = 1084.25
$ws.Range("I122").Value = 1167.5555
$ws.Range("J122").Value = 334.5
$ws.Range("K122").Value = 3502.6665
$ws.Range("L122").Value = 1003.5
$ws.Range("M122").Value = -1052.6665
$ws.Range("N122").Value = -5903.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 381.81818
$ws.Range("I5").Value = 391.375
$ws.Range("J5").Value = 356.33334
$ws.Range("K5").Value = 391.375
$ws.Range("L5").Value = 356.33334
$ws.Range("M5").Value = -279.375
$ws.Range("N5").Value = -580.33334
$ws.Range("H32").Value = 5388.174
$ws.Range("I32").Value = 5388.174
$ws.Range("K32").Value = 5388.174
$ws.Range("M32").Value = -5101.174
$ws.Range("H45").Value = 4870.8335
$ws.Range("I45").Value = 4642.2
$ws.Range("J45").Value = 6014
$ws.Range("K45").Value = 4642.2
$ws.Range("L45").Value = 6014
$ws.Range("M45").Value = -4265.2
$ws.Range("N45").Value = -6768
$ws.Range("H82").Value = 37000
$ws.Range("J82").Value = 37000
$ws.Range("L82").Value = 37000
$ws.Range("N82").Value = -37722
$ws.Range("H85").Value = 37000
$ws.Range("J85").Value = 37000
$ws.Range("L85").Value = 37000
$ws.Range("N85").Value = -39496
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840
$ws.Range("H132").Value = 8374.5
$ws.Range("I132").Value = 8374.5
$ws.Range("K132").Value = 25123.5
$ws.Range("M132").Value = -22593.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 381.81818
$ws.Range("I4").Value = 391.375
$ws.Range("J4").Value = 356.33334
$ws.Range("K4").Value = 391.375
$ws.Range("L4").Value = 356.33334
$ws.Range("M4").Value = -276.375
$ws.Range("N4").Value = -586.33334
$ws.Range("H6").Value = 9994.6
$ws.Range("J6").Value = 9994.6
$ws.Range("L6").Value = 9994.6
$ws.Range("N6").Value = -10220.6
$ws.Range("H19").Value = 16500
$ws.Range("H52").Value = 39989.168
$ws.Range("J52").Value = 39989.168
$ws.Range("L52").Value = 39989.168
$ws.Range("N52").Value = -40515.168
$ws.Range("H121").Value = 39989.168
$ws.Range("J121").Value = 39989.168
$ws.Range("L121").Value = 39989.168
$ws.Range("N121").Value = -43483.168
$ws.Range("H134").Value = 2063.5833
$ws.Range("I134").Value = 2074.0908
$ws.Range("K134").Value = 6222.2724
$ws.Range("M134").Value = -3687.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5875
$ws.Range("I31").Value = 10600
$ws.Range("J31").Value = 4693.75
$ws.Range("K31").Value = 10600
$ws.Range("L31").Value = 4693.75
$ws.Range("M31").Value = -10305
$ws.Range("N31").Value = -5283.75
$ws.Range("H34").Value = 5875
$ws.Range("I34").Value = 10600
$ws.Range("J34").Value = 4693.75
$ws.Range("K34").Value = 10600
$ws.Range("L34").Value = 4693.75
$ws.Range("M34").Value = -10398
$ws.Range("N34").Value = -5097.75
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H132").Value = 3086.8572
$ws.Range("I132").Value = 3135.3333
$ws.Range("K132").Value = 9405.999899999999
$ws.Range("M132").Value = -6875.999899999999
$ws.Range("H134").Value = 2904.4443
$ws.Range("I134").Value = 3205
$ws.Range("K134").Value = 9615
$ws.Range("M134").Value = -7080
$ws.Range("H141").Value = 496484.66
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3059.8
$ws.Range("J52").Value = 3059.8
$ws.Range("L52").Value = 9179.400000000001
$ws.Range("N52").Value = -9711.400000000001
$ws.Range("H54").Value = 1700
$ws.Range("J54").Value = 1700
$ws.Range("L54").Value = 5100
$ws.Range("N54").Value = -6218
$ws.Range("H114").Value = 2672.2856
$ws.Range("I114").Value = 3901.6667
$ws.Range("J114").Value = 1750.25
$ws.Range("K114").Value = 11705.0001
$ws.Range("L114").Value = 5250.75
$ws.Range("M114").Value = -8451.000100000001
$ws.Range("N114").Value = -11758.75
$ws.Range("H121").Value = 72143830
$ws.Range("I121").Value = 1061.2858
$ws.Range("J121").Value = 144286610
$ws.Range("K121").Value = 3183.8574
$ws.Range("L121").Value = 432859830
$ws.Range("M121").Value = -1873.8574
$ws.Range("N121").Value = -432862450
$ws.Range("H134").Value = 545
$ws.Range("I134").Value = 545
$ws.Range("K134").Value = 1635
$ws.Range("M134").Value = 3435
$ws.Range("H136").Value = 5996.4
$ws.Range("I136").Value = 5996.4
$ws.Range("K136").Value = 17989.2
$ws.Range("M136").Value = -12889.2
$ws.Range("H137").Value = 1259.8
$ws.Range("J137").Value = 2500
$ws.Range("L137").Value = 7500
$ws.Range("N137").Value = -17700
$ws.Range("H140").Value = 836947.25
$ws.Range("I140").Value = 912851.5600000001
$ws.Range("K140").Value = 2738554.68
$ws.Range("M140").Value = -2733374.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 22111.5
$ws.Range("J92").Value = 22111.5
$ws.Range("L92").Value = 22111.5
$ws.Range("N92").Value = -25855.5
$ws.Range("H102").Value = 1166.6666
$ws.Range("I102").Value = 1166.6666
$ws.Range("K102").Value = 1166.6666
$ws.Range("M102").Value = 455.3334
$ws.Range("H122").Value = 4147.5
$ws.Range("I122").Value = 3399.2
$ws.Range("J122").Value = 4895.8
$ws.Range("K122").Value = 10197.6
$ws.Range("L122").Value = 14687.4
$ws.Range("M122").Value = -7747.599999999999
$ws.Range("N122").Value = -19587.4
$ws.Range("H132").Value = 2991.7058
$ws.Range("I132").Value = 2825.75
$ws.Range("J132").Value = 3390
$ws.Range("K132").Value = 8477.25
$ws.Range("L132").Value = 10170
$ws.Range("M132").Value = -5947.25
$ws.Range("N132").Value = -15230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2889.25
$ws.Range("I40").Value = 2016.2858
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 2016.2858
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -1880.2858
$ws.Range("N40").Value = -9272
$ws.Range("H46").Value = 7200
$ws.Range("I46").Value = 7200
$ws.Range("K46").Value = 7200
$ws.Range("M46").Value = -7012
$ws.Range("H55").Value = 839.2
$ws.Range("I55").Value = 422.5
$ws.Range("J55").Value = 1117
$ws.Range("K55").Value = 422.5
$ws.Range("L55").Value = 1117
$ws.Range("M55").Value = -249.5
$ws.Range("N55").Value = -1463
$ws.Range("H61").Value = 1766
$ws.Range("I61").Value = 1481.6364
$ws.Range("K61").Value = 1481.6364
$ws.Range("M61").Value = -1279.6364
$ws.Range("H98").Value = 35927.223
$ws.Range("J98").Value = 35927.223
$ws.Range("L98").Value = 35927.223
$ws.Range("N98").Value = -41917.223
$ws.Range("H113").Value = 1766
$ws.Range("I113").Value = 1481.6364
$ws.Range("K113").Value = 1481.6364
$ws.Range("M113").Value = 688.3635999999999
$ws.Range("H118").Value = 44999.332
$ws.Range("I118").Value = 45000
$ws.Range("J118").Value = 44999
$ws.Range("K118").Value = 45000
$ws.Range("L118").Value = 44999
$ws.Range("M118").Value = -43343
$ws.Range("N118").Value = -48313
$ws.Range("H122").Value = 3553.6296
$ws.Range("I122").Value = 3798.4
$ws.Range("J122").Value = 3498
$ws.Range("K122").Value = 11395.2
$ws.Range("L122").Value = 10494
$ws.Range("M122").Value = -8945.200000000001
$ws.Range("N122").Value = -15394

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2593.75
$ws.Range("I122").Value = 2677.625
$ws.Range("K122").Value = 8032.875
$ws.Range("M122").Value = -5582.875

